$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Extend the "ss3" (full-feedback) scenario block with a new set of runs
# (rows 44-49), reusing the formatting of the previous pair of rows (41:42)
# so the alternating banded-fill styling continues seamlessly.
# ---------------------------------------------------------------------------
$ws.Range("A41:H42").Copy() | Out-Null
$ws.Range("A44:H45").PasteSpecial(-4122) | Out-Null
$ws.Range("A41:H42").Copy() | Out-Null
$ws.Range("A46:H47").PasteSpecial(-4122) | Out-Null
$ws.Range("A41:H42").Copy() | Out-Null
$ws.Range("A48:H49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 44
$ws.Range("A44").Value2 = 31
$ws.Range("B44").Value2 = "var"
$ws.Range("C44").Value2 = "med"
$ws.Range("D44").Value2 = "naq"
$ws.Range("E44").Value2 = "ss3"
$ws.Range("F44").Value2 = 6
$ws.Range("G44").Formula = '=CONCATENATE("ASS",E44,"_HCR",F44,"_REC",C44,"_INN",B44,"_OER",D44)'
$ws.Range("H44").Value2 = "Introduce observation error + SS3 assessment (full-feedback)"

# Row 45
$ws.Range("A45").Value2 = 32
$ws.Range("B45").Value2 = "var"
$ws.Range("C45").Value2 = "mix"
$ws.Range("D45").Value2 = "naq"
$ws.Range("E45").Value2 = "ss3"
$ws.Range("F45").Value2 = 6
$ws.Range("G45").Formula = '=CONCATENATE("ASS",E45,"_HCR",F45,"_REC",C45,"_INN",B45,"_OER",D45)'

# Row 46
$ws.Range("A46").Value2 = 33
$ws.Range("B46").Value2 = "var"
$ws.Range("C46").Value2 = "med"
$ws.Range("D46").Value2 = "naq"
$ws.Range("E46").Value2 = "ss3"
$ws.Range("F46").Value2 = 5
$ws.Range("G46").Formula = '=CONCATENATE("ASS",E46,"_HCR",F46,"_REC",C46,"_INN",B46,"_OER",D46)'

# Row 47
$ws.Range("A47").Value2 = 34
$ws.Range("B47").Value2 = "var"
$ws.Range("C47").Value2 = "mix"
$ws.Range("D47").Value2 = "naq"
$ws.Range("E47").Value2 = "ss3"
$ws.Range("F47").Value2 = 5
$ws.Range("G47").Formula = '=CONCATENATE("ASS",E47,"_HCR",F47,"_REC",C47,"_INN",B47,"_OER",D47)'

# Row 48
$ws.Range("A48").Value2 = 35
$ws.Range("B48").Value2 = "var"
$ws.Range("C48").Value2 = "low"
$ws.Range("D48").Value2 = "naq"
$ws.Range("E48").Value2 = "ss3"
$ws.Range("F48").Value2 = 6
$ws.Range("G48").Formula = '=CONCATENATE("ASS",E48,"_HCR",F48,"_REC",C48,"_INN",B48,"_OER",D48)'

# Row 49
$ws.Range("A49").Value2 = 36
$ws.Range("B49").Value2 = "var"
$ws.Range("C49").Value2 = "low"
$ws.Range("D49").Value2 = "naq"
$ws.Range("E49").Value2 = "ss3"
$ws.Range("F49").Value2 = 5
$ws.Range("G49").Formula = '=CONCATENATE("ASS",E49,"_HCR",F49,"_REC",C49,"_INN",B49,"_OER",D49)'

# Merge the comment column across the whole new block, mirroring the other blocks
$ws.Range("H44:H49").Merge() | Out-Null

# Keep the view pointed at the freshly added rows, like the source edit
$ws.Range("K41").Select() | Out-Null

Write-Output "done"
